$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 was the most-recently-added payment and had its phone number
# stored as text; now that a newer payment (row 24) is being recorded it
# gets normalized to a plain number, matching every earlier row.
$ws.Cells.Item(23, 1).Value = 71277628

# Start row 24 as a copy of row 23 so it inherits the same "touched but
# blank" shape for the optional columns (B: amount, F: discount_applied)
# before we overwrite the columns that actually differ for this payment.
$ws.Range("A23:I23").Copy($ws.Range("A24:I24"))

# New payment: phone 71277628, Cash, 2025-08-18T16:53:54, original_amount
# 766, no discount/birthday discount, final_amount 766, no points redeemed.
# Column A keeps the phone number as text for the newest row, same as row
# 23 did before it was normalized above.
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "71277628"
$ws.Cells.Item(24, 1).ClearFormats()

$ws.Cells.Item(24, 3).Value = "Cash"
$ws.Cells.Item(24, 4).Value = "2025-08-18T16:53:54"
$ws.Cells.Item(24, 5).Value = 766
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 766
